$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the cells we touch so that strings such as
# "1.011", "18.20" or "0.3630" keep their exact literal representation
# instead of being reinterpreted as numbers (which would drop trailing zeros).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('B25').NumberFormat = '@'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('B26').NumberFormat = '@'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('E51').NumberFormat = '@'

$ws.Range('D2').Value = '27.115.59'
$ws.Range('E2').Value = '  -0.45%  '
$ws.Range('D3').Value = '1.821.96'
$ws.Range('E3').Value = '  -0.71%  '
$ws.Range('D4').Value = '1.011'
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '312.13'
$ws.Range('E6').Value = '  -0.28%  '
$ws.Range('D7').Value = '0.4626'
$ws.Range('E7').Value = '  -1.84%  '
$ws.Range('D8').Value = '0.3630'
$ws.Range('E8').Value = '  -1.65%  '
$ws.Range('D9').Value = '0.07296'
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').Value = '0.8704'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').Value = '20.08'
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').Value = '1.878.00'
$ws.Range('E12').Value = '  +2.30%  '
$ws.Range('D13').Value = '0.07633'
$ws.Range('E13').Value = '  +4.12%  '
$ws.Range('D14').Value = '5.339'
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('D15').Value = '92.41'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '6.476'
$ws.Range('E16').Value = '  -1.54%  '
$ws.Range('E17').Value = '  -0.45%  '
$ws.Range('D18').Value = '0.000008648'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('E19').Value = '  -0.23%  '
$ws.Range('D20').Value = '27.408.67'
$ws.Range('E20').Value = '  +0.57%  '
$ws.Range('D21').Value = '14.47'
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('D22').Value = '5.214'
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').Value = '10.55'
$ws.Range('E23').Value = '  -1.40%  '
$ws.Range('D24').Value = '2.095.93'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '1.873'
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '151.45'
$ws.Range('E26').Value = '  -1.11%  '
$ws.Range('D27').Value = '18.20'
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('D28').Value = '2.074'
$ws.Range('E28').Value = '  -4.37%  '
$ws.Range('D29').Value = '5.097'
$ws.Range('D30').Value = '116.12'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('D31').Value = '0.08909'
$ws.Range('E31').Value = '  -0.19%  '
$ws.Range('D32').Value = '2.958'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '0.7354'
$ws.Range('E33').Value = '  -3.20%  '
$ws.Range('D34').Value = '4.454'
$ws.Range('E34').Value = '  -2.17%  '
$ws.Range('D35').Value = '1.137'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('E37').Value = '  +2.54%  '
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('D39').Value = '0.05248'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('D40').Value = '0.01913'
$ws.Range('E40').Value = '  -2.41%  '
$ws.Range('D41').Value = '2.923'
$ws.Range('E41').Value = '  -2.61%  '
$ws.Range('D42').Value = '7.149'
$ws.Range('E42').Value = '  -2.51%  '
$ws.Range('D43').Value = '0.5198'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').Value = '0.1627'
$ws.Range('E44').Value = '  -2.26%  '
$ws.Range('D45').Value = '8.273'
$ws.Range('E45').Value = '  -3.27%  '
$ws.Range('D46').Value = '0.4835'
$ws.Range('E46').Value = '  -2.43%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').Value = '10.16'
$ws.Range('E48').Value = '  -3.24%  '
$ws.Range('E49').Value = '  -0.58%  '
$ws.Range('D50').Value = '1.634'
$ws.Range('E50').Value = '  -2.43%  '
$ws.Range('D51').Value = '0.06262'
$ws.Range('E51').Value = '  -0.97%  '
